$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the attendance dates as plain text (DD/MM/YYYY). The
# tutorial's updated solution re-writes them using dashes (DD-MM-YYYY).
# Several of these (e.g. "01-08-2022") are day/month-ambiguous, so Excel's
# normal auto-detection would silently reinterpret them as real dates and
# store a date serial instead of text. Force "Text" number format before
# assigning, then drop back to the default "Normal" style afterward so the
# cell formatting matches the original (unstyled) cells exactly.
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    $cell.Style = "Normal"
}

# Attendance-count corrections that accompanied the date reformat.
$ws.Cells.Item(3, 4).Value = 3   # D3  Total Attendance Count
$ws.Cells.Item(3, 7).Value = 3   # G3  Invalid

$ws.Cells.Item(4, 4).Value = 1   # D4  Total Attendance Count
$ws.Cells.Item(4, 5).Value = 1   # E4  Real
$ws.Cells.Item(4, 8).Value = 0   # H4  Absent

$ws.Cells.Item(10, 4).Value = 1  # D10 Total Attendance Count
$ws.Cells.Item(10, 5).Value = 1  # E10 Real
$ws.Cells.Item(10, 8).Value = 0  # H10 Absent
